$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the time-slot labels in column C (rows 6 and 7)
$ws.Range("C6").Value = "19:45-19:50"
$ws.Range("C7").Value = "19:50-19:55"

# Move the active-cell selection from C11 to C10
$ws.Range("C10").Select()
